$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 100000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 100000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 100000
$ws.Range("N21").Value = -100936
$ws.Range("M21").Value = ""
$ws.Range("H23").Value = 100000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 100000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 100000
$ws.Range("N23").Value = -100468
$ws.Range("M23").Value = ""
$ws.Range("H69").Value = 6392
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 6392
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 19176
$ws.Range("N69").Value = -20924
$ws.Range("M69").Value = ""
$ws.Range("H72").Value = 6392
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 6392
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 57528
$ws.Range("N72").Value = -66264
$ws.Range("M72").Value = ""
$ws.Range("H74").Value = 3178.375
$ws.Range("I74").Value = 2864
$ws.Range("J74").Value = 3870
$ws.Range("K74").Value = 2864
$ws.Range("L74").Value = 3870
$ws.Range("M74").Value = -1928
$ws.Range("N74").Value = -5742
$ws.Range("H76").Value = 5490.6
$ws.Range("J76").Value = 5700.5713
$ws.Range("L76").Value = 5700.5713
$ws.Range("N76").Value = -6330.5713
$ws.Range("H77").Value = 3178.375
$ws.Range("I77").Value = 2864
$ws.Range("J77").Value = 3870
$ws.Range("K77").Value = 14320
$ws.Range("L77").Value = 19350
$ws.Range("M77").Value = -9640
$ws.Range("N77").Value = -28710
$ws.Range("H79").Value = 5490.6
$ws.Range("J79").Value = 5700.5713
$ws.Range("L79").Value = 5700.5713
$ws.Range("N79").Value = -7884.5713
$ws.Range("H92").Value = 1515.75
$ws.Range("I92").Value = 1596.6923
$ws.Range("J92").Value = 1165
$ws.Range("K92").Value = 1596.6923
$ws.Range("L92").Value = 1165
$ws.Range("M92").Value = -348.6922999999999
$ws.Range("N92").Value = -3661
$ws.Range("H129").Value = 2940.9185
$ws.Range("J129").Value = 945.3570999999999
$ws.Range("L129").Value = 2836.0713
$ws.Range("N129").Value = -12836.0713
$ws.Range("H132").Value = 3973148.2
$ws.Range("I132").Value = 4907277
$ws.Range("J132").Value = 3100.25
$ws.Range("K132").Value = 14721831
$ws.Range("L132").Value = 9300.75
$ws.Range("M132").Value = -14719301
$ws.Range("N132").Value = -14360.75
$ws.Range("H138").Value = 3148.923
$ws.Range("I138").Value = 2382.5908
$ws.Range("J138").Value = 3393.261
$ws.Range("K138").Value = 7147.7724
$ws.Range("L138").Value = 10179.783
$ws.Range("M138").Value = -2007.7724
$ws.Range("N138").Value = -20459.783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38889.03
$ws.Range("I32").Value = 16705.55
$ws.Range("J32").Value = 81143.28999999999
$ws.Range("K32").Value = 16705.55
$ws.Range("L32").Value = 81143.28999999999
$ws.Range("M32").Value = -16418.55
$ws.Range("N32").Value = -81717.28999999999
$ws.Range("H61").Value = 2180.0527
$ws.Range("I61").Value = 2351.4285
$ws.Range("J61").Value = 2141.3547
$ws.Range("K61").Value = 2351.4285
$ws.Range("L61").Value = 2141.3547
$ws.Range("M61").Value = -2139.4285
$ws.Range("N61").Value = -2565.3547
$ws.Range("H132").Value = 15356.048
$ws.Range("I132").Value = 19420.562
$ws.Range("K132").Value = 58261.686
$ws.Range("M132").Value = -55731.686
$ws.Range("H136").Value = 2180.0527
$ws.Range("I136").Value = 2351.4285
$ws.Range("J136").Value = 2141.3547
$ws.Range("K136").Value = 7054.2855
$ws.Range("L136").Value = 6424.0641
$ws.Range("M136").Value = -4504.2855
$ws.Range("N136").Value = -11524.0641

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 142858610
$ws.Range("I107").Value = 200001420
$ws.Range("K107").Value = 200001420
$ws.Range("M107").Value = -199999500
$ws.Range("H134").Value = 5475.7646
$ws.Range("I134").Value = 6699.3184
$ws.Range("J134").Value = 3232.5833
$ws.Range("K134").Value = 20097.9552
$ws.Range("L134").Value = 9697.749899999999
$ws.Range("M134").Value = -17562.9552
$ws.Range("N134").Value = -14767.7499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 560
$ws.Range("I22").Value = 194.5
$ws.Range("J22").Value = 925.5
$ws.Range("K22").Value = 194.5
$ws.Range("L22").Value = 925.5
$ws.Range("M22").Value = 155.5
$ws.Range("N22").Value = -1625.5
$ws.Range("H31").Value = 19066.164
$ws.Range("I31").Value = 947.2222
$ws.Range("J31").Value = 50628.195
$ws.Range("K31").Value = 947.2222
$ws.Range("L31").Value = 50628.195
$ws.Range("M31").Value = -652.2222
$ws.Range("N31").Value = -51218.195
$ws.Range("H33").Value = 4304.4287
$ws.Range("I33").Value = 3355.1667
$ws.Range("J33").Value = 10000
$ws.Range("K33").Value = 3355.1667
$ws.Range("L33").Value = 10000
$ws.Range("M33").Value = -2976.1667
$ws.Range("N33").Value = -10758
$ws.Range("H34").Value = 19066.164
$ws.Range("I34").Value = 947.2222
$ws.Range("J34").Value = 50628.195
$ws.Range("K34").Value = 947.2222
$ws.Range("L34").Value = 50628.195
$ws.Range("M34").Value = -745.2222
$ws.Range("N34").Value = -51032.195

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 642807.7
$ws.Range("J131").Value = 726560.9
$ws.Range("L131").Value = 2179682.7
$ws.Range("N131").Value = -2189762.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 10740
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 10740
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 10740
$ws.Range("N52").Value = -11258
$ws.Range("M52").Value = ""
$ws.Range("H80").Value = 111359340
$ws.Range("I80").Value = 143173090
$ws.Range("J80").Value = 11250
$ws.Range("K80").Value = 143173090
$ws.Range("L80").Value = 11250
$ws.Range("M80").Value = -143172092
$ws.Range("N80").Value = -13246
$ws.Range("H83").Value = 111359340
$ws.Range("I83").Value = 143173090
$ws.Range("J83").Value = 11250
$ws.Range("K83").Value = 715865450
$ws.Range("L83").Value = 56250
$ws.Range("M83").Value = -715860458
$ws.Range("N83").Value = -66234
$ws.Range("H97").Value = 37038370
$ws.Range("I97").Value = 47620530
$ws.Range("J97").Value = 803.6667
$ws.Range("K97").Value = 47620530
$ws.Range("L97").Value = 803.6667
$ws.Range("M97").Value = -47620034
$ws.Range("N97").Value = -1795.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1000
$ws.Range("I9").Value = 1000
$ws.Range("K9").Value = 1000
$ws.Range("M9").Value = -776
$ws.Range("N9").Value = ""
$ws.Range("H22").Value = 1532.6842
$ws.Range("I22").Value = 1669.8572
$ws.Range("J22").Value = 1452.6666
$ws.Range("K22").Value = 1669.8572
$ws.Range("L22").Value = 1452.6666
$ws.Range("M22").Value = -1374.8572
$ws.Range("N22").Value = -2042.6666
$ws.Range("H27").Value = 1532.6842
$ws.Range("I27").Value = 1669.8572
$ws.Range("J27").Value = 1452.6666
$ws.Range("K27").Value = 1669.8572
$ws.Range("L27").Value = 1452.6666
$ws.Range("M27").Value = -1562.8572
$ws.Range("N27").Value = -1666.6666
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = ""
$ws.Range("N35").Value = ""
$ws.Range("H100").Value = 2300
$ws.Range("J100").Value = 2350
$ws.Range("L100").Value = 2350
$ws.Range("N100").Value = -3432
$ws.Range("H132").Value = 4851
$ws.Range("I132").Value = 10465.667
$ws.Range("J132").Value = 2745.5
$ws.Range("K132").Value = 31397.001
$ws.Range("L132").Value = 8236.5
$ws.Range("M132").Value = -28867.001
$ws.Range("N132").Value = -13296.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6807.8
$ws.Range("J54").Value = 6914.75
$ws.Range("L54").Value = 6914.75
$ws.Range("N54").Value = -7954.75
$ws.Range("H81").Value = 251960.12
$ws.Range("I81").Value = 200960.4
$ws.Range("J81").Value = 336959.66
$ws.Range("K81").Value = 401920.8
$ws.Range("L81").Value = 673919.3199999999
$ws.Range("M81").Value = -400859.8
$ws.Range("N81").Value = -676041.3199999999
$ws.Range("H84").Value = 251960.12
$ws.Range("I84").Value = 200960.4
$ws.Range("J84").Value = 336959.66
$ws.Range("K84").Value = 2009604
$ws.Range("L84").Value = 3369596.6
$ws.Range("M84").Value = -2004300
$ws.Range("N84").Value = -3380204.6
$ws.Range("H107").Value = 91413.91
$ws.Range("I107").Value = 432
$ws.Range("J107").Value = 167232.17
$ws.Range("K107").Value = 1296
$ws.Range("L107").Value = 501696.51
$ws.Range("M107").Value = 624
$ws.Range("N107").Value = -505536.51
$ws.Range("H136").Value = 17784.607
$ws.Range("I136").Value = 29939.734
$ws.Range("K136").Value = 89819.202
$ws.Range("M136").Value = -87269.202

Write-Output "edits applied"